$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: header "25_03_2024" and values for each recepcionista
$ws.Range("F1").Value = "25_03_2024"
$ws.Range("F2").Value = 1580
$ws.Range("F3").Value = 1619
$ws.Range("F4").Value = 1474
$ws.Range("F5").Value = 344

# Update selection as recorded in the file (active cell F6)
$ws.Range("F6").Select()
